$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.027.46"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.31%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.124.21"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.83%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.25"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.12%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.36"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -5.37%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.119.51"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.88%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.83%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.98%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.25"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.37%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.77%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -5.00%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.17"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.64%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.638.16"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.85%  "

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.05%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.070.78"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.03%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.127.34"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.47%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.67"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -3.04%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "473.27"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.66%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.10"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.81%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.66"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.97"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.92%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.84"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -4.64%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.04%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.51%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.89"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -6.40%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.93"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.90%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.08"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.60%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.04%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.70"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.31%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.107"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -6.31%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.52"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -5.13%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.15%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.23%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.78%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0701"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -8.97%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0387"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.14%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "420.13"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -6.29%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.32%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.898.11"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.98%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -11.84%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.113"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -5.70%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.265"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.06%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.04%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -5.85%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.44"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.03%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.62%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -6.93%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "119.52"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.04%  "

